# mock_user_update_form.xlsx fixes:
#   1. Typo fix: the "general terms" consent text was missing a trailing
#      space after the final period.
#   2. Selection/scroll position nudged one column to the right (current
#      selection moves from F15 to J1).
#   3. Header/footer print margins normalized to 1.3cm (metric page setup).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the typo: add the missing trailing space to the consent text.
$ws.Range("J1").Value = "By clicking yes below, you agree with these general terms. "

# 2) Update the active selection / scroll position.
$ws.Range("J1").Select()

# 3) Header/footer margins -> 1.3 cm (matches the metric page-setup default).
$ws.PageSetup.HeaderMargin = $excel.CentimetersToPoints(1.3)
$ws.PageSetup.FooterMargin = $excel.CentimetersToPoints(1.3)
